$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.06"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").Value = "'23.74"
$ws.Range("D3").ClearFormats()
$ws.Range("D4").Value = "'5.325"
$ws.Range("D4").ClearFormats()
$ws.Range("D5").Value = "'0.05778"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").Value = "'6.475"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").Value = "'3.343"
$ws.Range("D7").ClearFormats()
$ws.Range("D8").Value = "'0.8106"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").Value = "'0.8864"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").Value = "'0.1394"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").Value = "'0.07362"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "'0.03086"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").Value = "'0.03056"
$ws.Range("D13").ClearFormats()
$ws.Range("D14").Value = "'0.09335"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").Value = "'3.876"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").Value = "'0.001536"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").Value = "'0.04714"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").Value = "'0.0006059"
$ws.Range("D18").ClearFormats()
$ws.Range("D19").Value = "'0.006172"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").Value = "'0.001297"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").Value = "'0.00008815"
$ws.Range("D21").ClearFormats()
$ws.Range("D24").Value = "'0.3181"
$ws.Range("D24").ClearFormats()
$ws.Range("D27").Value = "'0.004607"
$ws.Range("D27").ClearFormats()
$ws.Range("D28").Value = "'0.0002353"
$ws.Range("D28").ClearFormats()
$ws.Range("D40").Value = "'0.03777"
$ws.Range("D40").ClearFormats()
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1053"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("D42").Value = "'0.002564"
$ws.Range("D42").ClearFormats()
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003187"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
$ws.Range("D44").Value = "'0.007611"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").Value = "'0.00005478"
$ws.Range("D45").ClearFormats()
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").ClearFormats()
$ws.Range("D47").Value = "'0.5508"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").Value = "'0.001845"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("D49").Value = "'0.00002103"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").Value = "'0.0002003"
$ws.Range("D50").ClearFormats()
